# Scene.xlsx — add two NPC columns ("ActorID" / "CanClone") to the XML-bound
# table on Sheet1, fill in their values for the three existing rows, and
# touch the page-setup / selection the way the authoring session left them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Grow the table by two columns -----------------------------------
# The table currently spans A1:I4 (9 columns). Add two more list columns;
# ListColumns.Add() appends right after the last column, so the first
# Add() lands in J, the second in K.
$colJ = $tbl.ListColumns.Add()
$colK = $tbl.ListColumns.Add()

# Name the new columns via the header cells (renaming a ListColumn through
# its .Name property doesn't retarget the header text in this host, but
# writing the header cell value does - and it keeps the ListColumn's name
# in sync too). ActorID is entered before CanClone so the shared-string
# table picks up the same ordering as the authored workbook.
$ws.Cells.Item(1, 11).Value = "ActorID"
$ws.Cells.Item(1, 10).Value = "CanClone"

# --- Fill in the data for the 3 existing data rows --------------------
# Row 2 -> CloneScene (clone = 1), Row 3 -> PioneerNoob, Row 4 -> RebellerNoob
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 0

$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0

$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0

# --- Column J width (new "CanClone" column), matches the other ~14-wide
# columns on the sheet.
$ws.Columns.Item(10).ColumnWidth = 13.29

# --- Page setup: author touched the print setup (A4 / portrait) -------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection left on K9 after the edits ------------------------------
$ws.Range("K9").Select() | Out-Null
